$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 172
$ws.Range("F3").Value = 943
$ws.Range("F4").Value = 1098
$ws.Range("F5").Value = 1556
$ws.Range("F6").Value = 342
$ws.Range("F7").Value = 701
$ws.Range("F8").Value = 12745
$ws.Range("F9").Value = 2229
$ws.Range("F11").Value = 280
$ws.Range("F13").Value = 48970
$ws.Range("G13").Value = "已售罄"
$ws.Range("F14").Value = 1265
$ws.Range("F15").Value = 261
$ws.Range("F16").Value = 288
$ws.Range("F17").Value = 827
$ws.Range("F18").Value = 686
$ws.Range("F19").Value = 334
$ws.Range("F20").Value = 2945
$ws.Range("F21").Value = 812
$ws.Range("F22").Value = 4756
$ws.Range("F23").Value = 4756
$ws.Range("F24").Value = 1190
$ws.Range("F25").Value = 906
$ws.Range("F28").Value = 23
$ws.Range("F29").Value = 8
$ws.Range("F30").Value = 1136
$ws.Range("F31").Value = 68
$ws.Range("F32").Value = 129
$ws.Range("F33").Value = 294
$ws.Range("F34").Value = 33
$ws.Range("F36").Value = 47
$ws.Range("F38").Value = 4561
$ws.Range("F40").Value = 4654
$ws.Range("F41").Value = 5607
$ws.Range("F42").Value = 107
$ws.Range("F43").Value = 135
$ws.Range("F44").Value = 99
$ws.Range("F45").Value = 185
$ws.Range("F46").Value = 382
$ws.Range("F47").Value = 89
$ws.Range("F48").Value = 56
$ws.Range("F49").Value = 4135
$ws.Range("F50").Value = 158

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F2").Value = 4174
$ws.Range("F4").Value = 74
$ws.Range("F5").Value = 113
$ws.Range("F8").Value = 9
$ws.Range("F10").Value = 112
$ws.Range("F12").Value = 1069
$ws.Range("F13").Value = 6
$ws.Range("F16").Value = 7
$ws.Range("F18").Value = 54
$ws.Range("F19").Value = 1
$ws.Range("F20").Value = 82
$ws.Range("F22").Value = 12

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 770
$ws.Range("F3").Value = 501
$ws.Range("F4").Value = 120
$ws.Range("F5").Value = 23

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 501
$ws.Range("F3").Value = 172
$ws.Range("F4").Value = 943
$ws.Range("F5").Value = 1098
$ws.Range("F6").Value = 343
$ws.Range("F7").Value = 701
$ws.Range("F8").Value = 12746
$ws.Range("F9").Value = 12746
$ws.Range("F10").Value = 2229
$ws.Range("F12").Value = 280
$ws.Range("F13").Value = 1265
$ws.Range("F14").Value = 288
$ws.Range("F15").Value = 827
$ws.Range("F16").Value = 686
$ws.Range("F17").Value = 334
$ws.Range("F18").Value = 2945
$ws.Range("F19").Value = 812
$ws.Range("F20").Value = 74
$ws.Range("F21").Value = 4756
$ws.Range("F22").Value = 4756
$ws.Range("F23").Value = 1190
$ws.Range("F24").Value = 23
$ws.Range("F25").Value = 113
$ws.Range("F26").Value = 906
$ws.Range("F27").Value = 47
$ws.Range("F28").Value = 8
$ws.Range("F29").Value = 1136
$ws.Range("F30").Value = 4
$ws.Range("F31").Value = 68
$ws.Range("F32").Value = 129
$ws.Range("F34").Value = 294
$ws.Range("F35").Value = 33
$ws.Range("F37").Value = 4561
$ws.Range("F38").Value = 4654
$ws.Range("F39").Value = 107
$ws.Range("F40").Value = 135
$ws.Range("F41").Value = 99
$ws.Range("F42").Value = 185
$ws.Range("F43").Value = 382
$ws.Range("F44").Value = 7
$ws.Range("F46").Value = 89
$ws.Range("F47").Value = 4135
$ws.Range("F49").Value = 12
$ws.Range("F51").Value = 158
